# Applies updated cryptos list values (prices / 1h volume %) per the Wed Jun 19 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '65.278.22'
$ws.Range('E2').Value = '  -0.24%  '

$ws.Range('D3').Value = '3.532.57'
$ws.Range('E3').Value = '  +2.91%  '

$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '595.63'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.30%  '

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '138.83'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +2.71%  '

$ws.Range('D7').Value = '3.530.82'
$ws.Range('E7').Value = '  +2.77%  '

$ws.Range('E9').Value = '  +0.95%  '

$ws.Range('E10').Value = '  +3.10%  '

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.17'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -3.81%  '

$ws.Range('E12').Value = '  +3.38%  '

$ws.Range('D13').Value = '4.131.57'
$ws.Range('E13').Value = '  +3.15%  '

$ws.Range('E14').Value = '  +3.86%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = '3.545.60'
$ws.Range('E15').Value = '  +2.64%  '

$ws.Range('B16').Value = 'Avalanche'
$ws.Range('C16').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '26.89'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +2.09%  '

$ws.Range('E17').Value = '  +1.42%  '

$ws.Range('D18').Value = '65.137.53'
$ws.Range('E18').Value = '  -0.20%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.17'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +2.11%  '

$ws.Range('E20').Value = '  +2.09%  '

$ws.Range('E21').Value = '  +3.77%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '395.37'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.07%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.570'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +4.93%  '

$ws.Range('E24').Value = '  +1.64%  '

$ws.Range('D25').Value = '3.679.41'
$ws.Range('E25').Value = '  +3.07%  '

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.06%  '

$ws.Range('E27').Value = '  +7.72%  '

$ws.Range('E28').Value = '  +8.46%  '

$ws.Range('E29').Value = '  -0.01%  '

$ws.Range('E30').Value = '  +0.82%  '

$ws.Range('E31').Value = '  +0.85%  '

$ws.Range('D32').Value = '3.554.66'
$ws.Range('E32').Value = '  +3.56%  '

$ws.Range('E33').Value = '  +0.01%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '23.84'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +5.66%  '

$ws.Range('E35').Value = '  +0.47%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.23'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.02%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '170.71'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.85%  '

$ws.Range('E38').Value = '  +2.02%  '

$ws.Range('E39').Value = '  +1.04%  '

$ws.Range('E40').Value = '  +1.53%  '

$ws.Range('E41').Value = '  +4.02%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.822'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.04%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '26.68'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  +21.48%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '42.64'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.04%  '

$ws.Range('E45').Value = '  +0.08%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '4.42'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.66%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.18'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +7.97%  '

$ws.Range('E48').Value = '  +3.18%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '6.81'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +4.29%  '

$ws.Range('D50').Value = '2.345.60'
$ws.Range('E50').Value = '  +6.77%  '

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.10'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.23%  '
